$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.095.38'
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("D3").Value = '2.456.71'
$ws.Range("E3").Value = '  +0.58%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '509.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '134.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.26%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").Value = '2.455.41'
$ws.Range("E9").Value = '  +0.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0982'
$ws.Range("D10").Style = "Normal"
$ws.Range("E11").Value = '  -0.42%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.324'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.64'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.97%  '
$ws.Range("D14").Value = '2.897.01'
$ws.Range("E14").Value = '  +0.74%  '
$ws.Range("D15").Value = '57.957.97'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.98'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000136'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.32%  '
$ws.Range("D18").Value = '2.422.79'
$ws.Range("E18").Value = '  -0.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.19%  '
$ws.Range("E20").Value = '  +0.76%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '315.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.12%  '
$ws.Range("E23").Value = '  -0.06%  '
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.13%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("E27").Value = '  +0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.381'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.60'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +5.39%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '171.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("D31").Value = '0.0₃0738'
$ws.Range("E31").Value = '  +0.91%  '
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.15'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.18%  '
$ws.Range("E34").Value = '  +2.09%  '
$ws.Range("E35").Value = '  +0.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E38").Value = '  +5.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.90'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.95%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.812'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '136.56'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +14.06%  '
$ws.Range("E44").Value = '  +1.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.94'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.576'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.76%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '256.22'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0919'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("E49").Value = '  +1.19%  '
$ws.Range("E50").Value = '  +2.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.27'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.19%  '
